$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: header + "combo" values (Profile + "-" + FHIR_Element) ---
$ws.Range("E1").Value = "combo"

for ($r = 2; $r -le 34; $r++) {
    $profile = $ws.Cells.Item($r, 3).Text
    $element = $ws.Cells.Item($r, 4).Text
    $ws.Cells.Item($r, 5).Value = $profile + "-" + $element
}

# --- Rename "An Interpreter Required Flag" -> "An Interpreter Needed Flag" ---
$ws.Range("B4").Value = "An Interpreter Needed Flag"
$ws.Range("B5").Value = "An Interpreter Needed Flag"

# --- Append three new rows (35-37) ---
$ws.Range("B35").Value = "Pregnancy Status"
$ws.Range("C35").Value = "US Core Observation Pregnancy Intent Profile"
$ws.Range("D35").Value = "Observation.performer"
$ws.Range("E35").Value = "US Core Observation Pregnancy Intent Profile-Observation.performer"

$ws.Range("B36").Value = "Sexual Orientation"
$ws.Range("C36").Value = "US Core Observation Sexual Orientation Profile"
$ws.Range("D36").Value = "Observation.performer"
$ws.Range("E36").Value = "US Core Observation Sexual Orientation Profile-Observation.performer"

$ws.Range("B37").Value = "Pregnancy Status"
$ws.Range("C37").Value = "US Core Observation Pregnancy Status Profile"
$ws.Range("D37").Value = "Observation.performer"
$ws.Range("E37").Value = "US Core Observation Pregnancy Status Profile-Observation.performer"

# Column A for the new rows must hold the literal text "True" (same as the
# existing A4/A5/etc. cells), not an Excel Boolean. Assigning the bare word
# "True" via .Value auto-converts to a Boolean, and an apostrophe-text-prefix
# leaves a stray quotePrefix style on the cell. Instead, stage a text-formula
# result in a scratch cell and paste-special its *value* into the target
# cells: that lands as a plain shared string with no extra styling, exactly
# like the pre-existing "True" cells.
$stage = $ws.Range("Z1")
$stage.Formula = "=""True"""
$stage.Copy()
$ws.Range("A35").PasteSpecial(-4163)
$ws.Range("A36").PasteSpecial(-4163)
$ws.Range("A37").PasteSpecial(-4163)
$stage.Clear()
